$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Sheet1" -> "Labels" and populate names ---
$wb.Worksheets.Item(1).Name = "Labels"
$wsLabels = $wb.Worksheets.Item("Labels")

$wsLabels.Range("A1").Value = "Alice"
$wsLabels.Range("A2").Value = "Bob"
$wsLabels.Range("A3").Value = "Charlie"
$wsLabels.Range("A4").Value = "Daniel"
$wsLabels.Range("A5").Value = "Eve"
$wsLabels.Range("A4").Select()

# --- Add "Formula" then "Sheet2" right after Labels so the final tab ---
# --- order becomes Labels, Sheet2, Formula. Re-fetch sheets by NAME ---
# --- after every structural change since worksheet variables can go ---
# --- stale (rebind to a tab position) once the sheet collection is ---
# --- reshuffled. ---
$wb.Worksheets.Add($null, $wb.Worksheets.Item("Labels")).Name = "Formula"
$wb.Worksheets.Add($null, $wb.Worksheets.Item("Labels")).Name = "Sheet2"

$wsSheet2 = $wb.Worksheets.Item("Sheet2")
$wsFormula = $wb.Worksheets.Item("Formula")

# --- Sheet2: two columns of numbers ---
for ($i = 1; $i -le 10; $i++) {
    $wsSheet2.Cells.Item($i, 1).Value = $i
    $wsSheet2.Cells.Item($i, 2).Value = $i + 10
}
$wsSheet2.Range("B11").Select()

# --- Formula sheet: labels + cross-sheet / same-sheet formulas ---
$wsFormula.Range("A1").Value = "Label Ref"
$wsFormula.Range("B1").Formula = "=Labels!A1"

$wsFormula.Range("A2").Value = "Label Ref Legacy"
$wsFormula.Range("B2").Formula = "=+Labels!A1"

$wsFormula.Range("A3").Value = "Label Ref Same Sheet"
$wsFormula.Range("B3").Formula = "=A3"

$wsFormula.Range("A4").Value = "Label Ref Same Sheet Legacy"
$wsFormula.Range("B4").Formula = "=A4"

# Widen column A to (roughly) fit the longest label
$wsFormula.Columns.Item(1).ColumnWidth = 21.5

$wsFormula.Range("B1").Select()

# Formula is the tab that should end up active/selected
$wsFormula.Activate()
